# New "prelievi" (withdrawal) entry dated 11 May 2018 (Excel serial 43231),
# appended as row 21 beneath the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the row above so the new row inherits the same per-column cell
# styles/number-formats (date format on A, text formats on B/C/D, plain
# number on E) instead of picking up a generic default style.
$ws.Rows(20).Copy()
$ws.Rows(21).Insert(-4121)

$ws.Cells.Item(21, 1).Value = 43231
$ws.Cells.Item(21, 2).Value = "Segreteria"
$ws.Cells.Item(21, 3).Value = "Scotch Magic"
$ws.Cells.Item(21, 4).Value = "N°."
$ws.Cells.Item(21, 5).Value = 5
